# Slide 7 ("Title 1" shape) currently holds a single run: "AI Safety Checklist".
# Split it into two runs -- "AI Safety " (unchanged formatting/color) and
# "Checklist" (same font/size, but recolored to #1E3A79) -- by grabbing the
# back half of the text via Characters() and changing its font color, which
# causes PowerPoint to break it out into its own run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

# "AI Safety Checklist"
#  123456789012345678 9
# "Checklist" starts at character 11 and is 9 characters long.
$tail = $tr.Characters(11, 9)

# COM Font.Color.RGB is packed as 0xBBGGRR (R + G*256 + B*65536), matching
# VBA's RGB() function, so build the target color 0x1E3A79 (R=0x1E,G=0x3A,B=0x79)
# that way instead of writing the raw hex value directly.
$r = 0x1E
$g = 0x3A
$b = 0x79
$tail.Font.Color.RGB = $r + ($g * 256) + ($b * 65536)
